# "Worked on temporal resolution"
#
# Changes applied to the Demand sheet:
#  - Update the demand value for timestep t=1
#  - Extend the demand time series from 1 timestep to 12 timesteps (t = 1..12),
#    all using the same (updated) annual demand value
#  - Auto-fit column B of the Demand sheet to its new (wider) content
#  - Make the "Demand" sheet the active/selected sheet (previously "SupIm" was active)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

$demandValue = 436619792

# Existing row 3 (t = 1) gets the new value
$ws.Cells.Item(3, 2).Value2 = $demandValue

# Add rows for t = 2 .. 12 (rows 4 .. 14), all sharing the same demand value
for ($t = 2; $t -le 12; $t++) {
    $row = $t + 2
    $ws.Cells.Item($row, 1).Value2 = $t
    $ws.Cells.Item($row, 2).Value2 = $demandValue
}

# Fit column B width to the new (wider) numbers
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# Make "Demand" the active sheet/tab (was "SupIm" before) and update its selection
$ws.Activate() | Out-Null
$ws.Range("E15").Select() | Out-Null
